$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F, matching style of existing header cells (e.g. E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "scenario"

# Fill F2:F101 with the scenario label "S4"
$lastRow = 101
$ws.Range("F2:F$lastRow").Value = "S4"
